$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the underlying input values; the SUM formulas in column E
# recalculate automatically.
$ws.Range("B9").Value = 4
$ws.Range("C12").Value = 6

# Move the active selection on the frozen-pane view to G3.
$ws.Range("G3").Select()
